$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helper = $ws.Cells.Item(1, 20)
$helper.Value = "'TRUE"
$helper.Copy()
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 6).PasteSpecial(-4163)
}
$helper.Clear()
$excel.CutCopyMode = 0

$ws.Range("F2:F17").Select()
